# Delete the rows that were removed from the "NEW" sheet of the
# interactive map export. Deleting from the bottom up so row numbers
# of not-yet-deleted rows remain valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

$rowsToDelete = @(77, 75, 69, 68, 18)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
